# Apply the dated "two-digit number divided by one-digit number"
# worksheet update: bump the header date by one day, and replace each
# of the 25 division problems in the table with the new problem.
#
# A positional (row/column) approach is used for the table cells rather
# than a global Find/Replace, because some of the new values coincide
# with old values used elsewhere in the table (e.g. "59÷2=" becomes
# "61÷5=", which was itself an original value elsewhere that becomes
# "65÷7="). A sequential text-based find/replace could cause one
# replacement to clobber another; addressing cells by position avoids
# that entirely.

$d = $word.ActiveDocument

# --- Update the header date line -------------------------------------
$d.Content.Find.Execute(
    "2025-12-26 Friday", $true, $false, $false, $false, $false,
    $true, 1, $false, "2025-12-27 Saturday", 2) | Out-Null

# --- Update the division problems in the table ------------------------
$t = $d.Tables.Item(1)

# Each inner array is one table row (1-indexed), holding the new text
# for columns 1..5 in order. Only the five rows that actually contain
# problems are listed (the rows in between are blank spacer rows).
$rowUpdates = @{
    1  = @("91÷7=", "65÷7=", "60÷5=", "91÷4=", "57÷3=")
    5  = @("32÷7=", "64÷8=", "17÷7=", "66÷2=", "64÷9=")
    9  = @("99÷8=", "55÷8=", "77÷4=", "29÷2=", "95÷2=")
    13 = @("78÷7=", "61÷5=", "34÷2=", "43÷5=", "49÷8=")
    17 = @("77÷5=", "80÷6=", "17÷3=", "82÷2=", "91÷7=")
}

foreach ($rowIndex in $rowUpdates.Keys) {
    $values = $rowUpdates[$rowIndex]
    for ($col = 1; $col -le 5; $col++) {
        $t.Cell($rowIndex, $col).Range.Text = $values[$col - 1]
    }
}
